# "Adding a Table" slide -> "Adding a Bullet Slide" slide
# (python-pptx "Adding a Bullet Slide" tutorial shape set, applied via COM)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1. Retitle the slide. Using a trailing paragraph mark + deleting the
#    spare trailing paragraph keeps the run free of an explicit <a:rPr>,
#    matching a from-scratch python-pptx style run.
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "Adding a Bullet Slide`r"
$titleTr.Paragraphs(2).Delete()

# 2. Drop the sample table ("Table 2") that used to sit below the title.
$s.Shapes.Item(2).Delete()

# 3. Add the bullet-slide content placeholder (idx=1, "Content Placeholder 2").
#    There's no direct "add placeholder by idx" verb, so borrow the
#    "Title and Content" layout's placeholder momentarily, then restore the
#    slide's original "Title Only" layout -- the inherited placeholder shape
#    (and its <p:ph idx="1"/>) stays on the slide either way.
$masterLayouts = $p.SlideMaster.CustomLayouts
$originalLayout = $s.CustomLayout
$titleAndContentLayout = $masterLayouts.Item(2)
$s.CustomLayout = $titleAndContentLayout
$s.CustomLayout = $originalLayout

$body = $s.Shapes.Item($s.Shapes.Count)
$bodyTr = $body.TextFrame.TextRange
$bodyTr.Text = "Find the bullet slide layout`rUse _TextFrame.add_paragraph() for subsequent bullets"
$bodyTr.Paragraphs(2).IndentLevel = 3
